$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# coinranking.com crypto price/volume refresh (GitHub Actions bot).
# A handful of the new Price values are plain decimals (e.g. "0.580",
# "3.40") that Excel would otherwise auto-convert to numbers on assignment,
# silently dropping the trailing zero the source feed renders as text.
# A leading apostrophe (Excel's own "treat as text" entry convention) keeps
# those values verbatim without touching the cell's number format.

$ws.Range('D2').Value = '51.750.14'
$ws.Range('E2').Value = '  +3.72%  '
$ws.Range('D3').Value = '2.754.48'
$ws.Range('E3').Value = '  +2.98%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'115.98"
$ws.Range('E5').Value = '  +2.10%  '
$ws.Range('D6').Value = "'333.87"
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = "'0.580"
$ws.Range('E9').Value = '  +4.80%  '
$ws.Range('D10').Value = "'41.64"
$ws.Range('E10').Value = '  +2.17%  '
$ws.Range('D11').Value = "'20.27"
$ws.Range('E11').Value = '  +0.78%  '
$ws.Range('D12').Value = "'0.0830"
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('E13').Value = '  +2.66%  '
$ws.Range('D14').Value = "'7.64"
$ws.Range('D15').Value = '3.181.89'
$ws.Range('E15').Value = '  +2.75%  '
$ws.Range('D16').Value = '2.735.99'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('E17').Value = '  +1.53%  '
$ws.Range('D18').Value = '51.676.34'
$ws.Range('E18').Value = '  +3.65%  '
$ws.Range('D19').Value = "'13.95"
$ws.Range('E19').Value = '  +5.84%  '
$ws.Range('B20').Value = 'ImmutableX'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D20').Value = "'2.99"
$ws.Range('E20').Value = '  +3.16%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = "'6.89"
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = "'279.14"
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = "'70.24"
$ws.Range('E24').Value = '  -2.28%  '
$ws.Range('D25').Value = "'2.67"
$ws.Range('E25').Value = '  +4.23%  '
$ws.Range('D26').Value = "'27.04"
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').Value = "'10.39"
$ws.Range('E28').Value = '  +2.19%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').Value = "'0.140"
$ws.Range('E31').Value = '  -1.17%  '
$ws.Range('D32').Value = "'50.43"
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('E33').Value = '  +3.13%  '
$ws.Range('D34').Value = "'0.0829"
$ws.Range('E34').Value = '  +2.60%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.25%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = "'2.11"
$ws.Range('E37').Value = '  +1.68%  '
$ws.Range('D38').Value = "'5.01"
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = "'3.25"
$ws.Range('E39').Value = '  +2.84%  '
$ws.Range('D40').Value = "'129.62"
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('D41').Value = "'23.84"
$ws.Range('E41').Value = '  +4.86%  '
$ws.Range('D42').Value = "'0.0350"
$ws.Range('E42').Value = '  +10.44%  '
$ws.Range('E43').Value = '  +3.99%  '
$ws.Range('D44').Value = "'0.113"
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +16.25%  '
$ws.Range('D46').Value = "'3.40"
$ws.Range('E46').Value = '  +2.66%  '
$ws.Range('D47').Value = '2.108.42'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('D48').Value = "'2.26"
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = "'5.62"
$ws.Range('E49').Value = '  +5.07%  '
$ws.Range('D50').Value = "'9.06"
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = "'1.55"
$ws.Range('E51').Value = '  +9.13%  '
